$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 107.333336
$ws.Range("I11").Value = 107.333336
$ws.Range("K11").Value = 107.333336
$ws.Range("M11").Value = 32.666664
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 8500
$ws.Range("J64").Value = 8500
$ws.Range("L64").Value = 8500
$ws.Range("N64").Value = -8996
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 8500
$ws.Range("J67").Value = 8500
$ws.Range("L67").Value = 8500
$ws.Range("N67").Value = -10216
# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 1117.8
$ws.Range("I100").Value = 1218.5714
$ws.Range("J100").Value = 882.6667
$ws.Range("K100").Value = 1218.5714
$ws.Range("L100").Value = 882.6667
$ws.Range("M100").Value = -677.5714
$ws.Range("N100").Value = -1964.6667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1010.53845
$ws.Range("I2").Value = 928.0833
$ws.Range("K2").Value = 928.0833
$ws.Range("M2").Value = -815.0833
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3298.4167
$ws.Range("I45").Value = 1891.2
$ws.Range("K45").Value = 1891.2
$ws.Range("M45").Value = -1514.2
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 877.44446
$ws.Range("I97").Value = 800
$ws.Range("K97").Value = 800
$ws.Range("M97").Value = -304
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1010.53845
$ws.Range("I116").Value = 928.0833
$ws.Range("K116").Value = 928.0833
$ws.Range("M116").Value = 1365.9167
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1010.53845
$ws.Range("I3").Value = 928.0833
$ws.Range("K3").Value = 928.0833
$ws.Range("M3").Value = -814.0833
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 765.375
$ws.Range("I99").Value = 765.375
$ws.Range("K99").Value = 765.375
$ws.Range("M99").Value = 732.625
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 5400.1665
$ws.Range("I107").Value = 3100.25
$ws.Range("K107").Value = 3100.25
$ws.Range("M107").Value = -1180.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2745
$ws.Range("I16").Value = 2740
$ws.Range("K16").Value = 2740
$ws.Range("M16").Value = -2453
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 62755
$ws.Range("J68").Value = 62755
$ws.Range("L68").Value = 62755
$ws.Range("N68").Value = -64253
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 62755
$ws.Range("J71").Value = 62755
$ws.Range("L71").Value = 188265
$ws.Range("N71").Value = -195753
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 3385.1667
$ws.Range("I105").Value = 3385.1667
$ws.Range("K105").Value = 3385.1667
$ws.Range("M105").Value = -1638.1667
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2745
$ws.Range("I113").Value = 2740
$ws.Range("K113").Value = 2740
$ws.Range("M113").Value = -570
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3194.75
$ws.Range("I132").Value = 3739
$ws.Range("J132").Value = 2650.5
$ws.Range("K132").Value = 11217
$ws.Range("L132").Value = 7951.5
$ws.Range("M132").Value = -8687
$ws.Range("N132").Value = -13011.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 16341
$ws.Range("J57").Value = 16341
$ws.Range("L57").Value = 16341
$ws.Range("N57").Value = -17981
# Row 68 (Leve Item ID 10659)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 8500
$ws.Range("I70").Value = 8500
$ws.Range("K70").Value = 8500
$ws.Range("M70").Value = -8230
# Row 71 (Leve Item ID 10659)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 8500
$ws.Range("I73").Value = 8500
$ws.Range("K73").Value = 8500
$ws.Range("M73").Value = -7564
# Row 87 (Leve Item ID 11894)
$ws.Range("H87").Value = 20354
$ws.Range("J87").Value = 20354
$ws.Range("L87").Value = 20354
$ws.Range("N87").Value = -22850
# Row 90 (Leve Item ID 11894)
$ws.Range("H90").Value = 20354
$ws.Range("J90").Value = 20354
$ws.Range("L90").Value = 61062
$ws.Range("N90").Value = -73542
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1948.3914
$ws.Range("I102").Value = 1486.15
$ws.Range("K102").Value = 1486.15
$ws.Range("M102").Value = 135.8499999999999
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 3338.5
$ws.Range("I113").Value = 1626.4286
$ws.Range("K113").Value = 1626.4286
$ws.Range("M113").Value = 543.5714
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 2232.8462
$ws.Range("J122").Value = 3958.5
$ws.Range("L122").Value = 11875.5
$ws.Range("N122").Value = -16775.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1072
$ws.Range("I22").Value = 734.5714
$ws.Range("J22").Value = 1662.5
$ws.Range("K22").Value = 734.5714
$ws.Range("L22").Value = 1662.5
$ws.Range("M22").Value = -439.5714
$ws.Range("N22").Value = -2252.5
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1072
$ws.Range("I27").Value = 734.5714
$ws.Range("J27").Value = 1662.5
$ws.Range("K27").Value = 734.5714
$ws.Range("L27").Value = 1662.5
$ws.Range("M27").Value = -627.5714
$ws.Range("N27").Value = -1876.5
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 6428.7144
$ws.Range("I40").Value = 4998.8
$ws.Range("K40").Value = 4998.8
$ws.Range("M40").Value = -4862.8
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 7047.125
$ws.Range("I61").Value = 6475.4
$ws.Range("K61").Value = 6475.4
$ws.Range("M61").Value = -6273.4
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3937.5
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 4583.3335
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 4583.3335
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -6081.3335
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3937.5
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 4583.3335
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 22916.6675
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -30404.6675
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 7047.125
$ws.Range("I113").Value = 6475.4
$ws.Range("K113").Value = 6475.4
$ws.Range("M113").Value = -4305.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 574.0909
$ws.Range("I107").Value = 574.0909
$ws.Range("K107").Value = 1722.2727
$ws.Range("M107").Value = 197.7273
